# Update team-specific transition probability matrix values
# (Bellarmine_B.xlsx) per the refreshed team-specific time data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2200647249190938
$ws.Range("C2").Value = 0.5469255663430421
$ws.Range("J2").Value = 0.003236245954692557
$ws.Range("P2").Value = 0.1682847896440129
$ws.Range("S2").Value = 0.06148867313915857
$ws.Range("B3").Value = 0.02777777777777778
$ws.Range("C3").Value = 0.03888888888888889
$ws.Range("J3").Value = 0.03888888888888889
$ws.Range("P3").Value = 0.7055555555555556
$ws.Range("S3").Value = 0.1888888888888889
$ws.Range("J4").Value = 0.04
$ws.Range("P4").Value = 0.8
$ws.Range("S4").Value = 0.16
$ws.Range("J5").Value = 0.1428571428571428
$ws.Range("P5").Value = 0.5714285714285714
$ws.Range("S5").Value = 0.2857142857142857
$ws.Range("B6").Value = 0.05150214592274678
$ws.Range("D6").Value = 0.01716738197424893
$ws.Range("E6").Value = 0.008583690987124463
$ws.Range("F6").Value = 0.04291845493562232
$ws.Range("J6").Value = 0.278969957081545
$ws.Range("O6").Value = 0.0128755364806867
$ws.Range("Q6").Value = 0.1373390557939914
$ws.Range("R6").Value = 0.09012875536480687
$ws.Range("S6").Value = 0.3605150214592275
$ws.Range("B7").Value = 0.1235294117647059
$ws.Range("D7").Value = 0.01764705882352941
$ws.Range("F7").Value = 0.02352941176470588
$ws.Range("J7").Value = 0.1588235294117647
$ws.Range("O7").Value = 0.01176470588235294
$ws.Range("Q7").Value = 0.1941176470588235
$ws.Range("R7").Value = 0.1058823529411765
$ws.Range("S7").Value = 0.3647058823529412
$ws.Range("B8").Value = 0.09145129224652088
$ws.Range("D8").Value = 0.01590457256461232
$ws.Range("E8").Value = 0.005964214711729622
$ws.Range("F8").Value = 0.06163021868787277
$ws.Range("J8").Value = 0.1053677932405567
$ws.Range("O8").Value = 0.02186878727634195
$ws.Range("Q8").Value = 0.1789264413518887
$ws.Range("R8").Value = 0.121272365805169
$ws.Range("S8").Value = 0.3976143141153082
$ws.Range("B9").Value = 0.06521739130434782
$ws.Range("D9").Value = 0.01304347826086956
$ws.Range("E9").Value = 0.004347826086956522
$ws.Range("F9").Value = 0.06521739130434782
$ws.Range("J9").Value = 0.1043478260869565
$ws.Range("O9").Value = 0.01739130434782609
$ws.Range("Q9").Value = 0.1826086956521739
$ws.Range("R9").Value = 0.1217391304347826
$ws.Range("S9").Value = 0.4260869565217391
$ws.Range("B10").Value = 0.09992810927390366
$ws.Range("D10").Value = 0.02300503235082674
$ws.Range("E10").Value = 0.0007189072609633358
$ws.Range("F10").Value = 0.07189072609633357
$ws.Range("J10").Value = 0.111430625449317
$ws.Range("O10").Value = 0.01653486700215672
$ws.Range("Q10").Value = 0.2070452911574407
$ws.Range("R10").Value = 0.1020848310567937
$ws.Range("S10").Value = 0.3673616103522646
$ws.Range("G11").Value = 0.1312217194570136
$ws.Range("J11").Value = 0.05882352941176471
$ws.Range("K11").Value = 0.1538461538461539
$ws.Range("L11").Value = 0.6470588235294118
$ws.Range("S11").Value = 0.009049773755656109
$ws.Range("G12").Value = 0.8125
$ws.Range("J12").Value = 0.1458333333333333
$ws.Range("K12").Value = 0.01388888888888889
$ws.Range("S12").Value = 0.02777777777777778
$ws.Range("G13").Value = 0.65
$ws.Range("J13").Value = 0.325
$ws.Range("S13").Value = 0.025
$ws.Range("F15").Value = 0.01271186440677966
$ws.Range("H15").Value = 0.1694915254237288
$ws.Range("I15").Value = 0.06779661016949153
$ws.Range("J15").Value = 0.3813559322033898
$ws.Range("K15").Value = 0.04661016949152542
$ws.Range("M15").Value = 0.008474576271186441
$ws.Range("O15").Value = 0.06779661016949153
$ws.Range("S15").Value = 0.2457627118644068
$ws.Range("F16").Value = 0.01415094339622642
$ws.Range("H16").Value = 0.2028301886792453
$ws.Range("I16").Value = 0.07547169811320754
$ws.Range("J16").Value = 0.4056603773584906
$ws.Range("K16").Value = 0.05660377358490566
$ws.Range("M16").Value = 0.02830188679245283
$ws.Range("N16").Value = 0.009433962264150943
$ws.Range("O16").Value = 0.06132075471698113
$ws.Range("S16").Value = 0.1462264150943396
$ws.Range("F17").Value = 0.025
$ws.Range("H17").Value = 0.1854166666666667
$ws.Range("I17").Value = 0.08749999999999999
$ws.Range("J17").Value = 0.4583333333333333
$ws.Range("K17").Value = 0.0625
$ws.Range("M17").Value = 0.02291666666666667
$ws.Range("O17").Value = 0.05416666666666667
$ws.Range("S17").Value = 0.1041666666666667
$ws.Range("F18").Value = 0.01865671641791045
$ws.Range("H18").Value = 0.2574626865671642
$ws.Range("I18").Value = 0.08582089552238806
$ws.Range("J18").Value = 0.4029850746268657
$ws.Range("K18").Value = 0.08955223880597014
$ws.Range("M18").Value = 0.007462686567164179
$ws.Range("O18").Value = 0.05970149253731343
$ws.Range("S18").Value = 0.07835820895522388
$ws.Range("F19").Value = 0.01482059282371295
$ws.Range("H19").Value = 0.2059282371294852
$ws.Range("I19").Value = 0.1060842433697348
$ws.Range("J19").Value = 0.406396255850234
$ws.Range("K19").Value = 0.08502340093603744
$ws.Range("M19").Value = 0.01638065522620905
$ws.Range("N19").Value = 0.0007800312012480499
$ws.Range("O19").Value = 0.0748829953198128
$ws.Range("S19").Value = 0.08970358814352575
